$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(485, 1).Value = 484
$ws.Cells.Item(485, 2).Value = "28.11.2023"
$ws.Cells.Item(485, 3).Value = "14:49:43"
$ws.Cells.Item(485, 4).Value = 34
$ws.Cells.Item(485, 5).Value = 21.4
$ws.Cells.Item(485, 6).Value = " "
$ws.Cells.Item(485, 7).Value = "   "

$ws.Cells.Item(486, 1).Value = 485
$ws.Cells.Item(486, 2).Value = "28.11.2023"
$ws.Cells.Item(486, 3).Value = "14:50:45"
$ws.Cells.Item(486, 4).Value = 33.9
$ws.Cells.Item(486, 5).Value = 21.2
$ws.Cells.Item(486, 6).Value = " "
$ws.Cells.Item(486, 7).Value = "   "

$ws.Cells.Item(487, 1).Value = 486
$ws.Cells.Item(487, 2).Value = "28.11.2023"
$ws.Cells.Item(487, 3).Value = "14:51:48"
$ws.Cells.Item(487, 4).Value = 33.8
$ws.Cells.Item(487, 5).Value = 21.3
$ws.Cells.Item(487, 6).Value = " "
$ws.Cells.Item(487, 7).Value = "   "

$ws.Cells.Item(488, 1).Value = 487
$ws.Cells.Item(488, 2).Value = "28.11.2023"
$ws.Cells.Item(488, 3).Value = "14:52:49"
$ws.Cells.Item(488, 4).Value = 33.7
$ws.Cells.Item(488, 5).Value = 21.4
$ws.Cells.Item(488, 6).Value = " "
$ws.Cells.Item(488, 7).Value = "   "

$ws.Cells.Item(489, 1).Value = 488
$ws.Cells.Item(489, 2).Value = "28.11.2023"
$ws.Cells.Item(489, 3).Value = "14:53:50"
$ws.Cells.Item(489, 4).Value = 33.8
$ws.Cells.Item(489, 5).Value = 21.4
$ws.Cells.Item(489, 6).Value = " "
$ws.Cells.Item(489, 7).Value = "   "

$ws.Cells.Item(490, 1).Value = 489
$ws.Cells.Item(490, 2).Value = "28.11.2023"
$ws.Cells.Item(490, 3).Value = "14:54:54"
$ws.Cells.Item(490, 4).Value = 33.9
$ws.Cells.Item(490, 5).Value = 21.4
$ws.Cells.Item(490, 6).Value = " "
$ws.Cells.Item(490, 7).Value = "   "

$ws.Cells.Item(491, 1).Value = 490
$ws.Cells.Item(491, 2).Value = "28.11.2023"
$ws.Cells.Item(491, 3).Value = "14:55:58"
$ws.Cells.Item(491, 4).Value = 33.9
$ws.Cells.Item(491, 5).Value = 21.4
$ws.Cells.Item(491, 6).Value = " "
$ws.Cells.Item(491, 7).Value = "   "

$ws.Cells.Item(492, 1).Value = 491
$ws.Cells.Item(492, 2).Value = "28.11.2023"
$ws.Cells.Item(492, 3).Value = "14:57:04"
$ws.Cells.Item(492, 4).Value = 33.9
$ws.Cells.Item(492, 5).Value = 21.4
$ws.Cells.Item(492, 6).Value = " "
$ws.Cells.Item(492, 7).Value = "   "

$ws.Cells.Item(493, 1).Value = 492
$ws.Cells.Item(493, 2).Value = "28.11.2023"
$ws.Cells.Item(493, 3).Value = "14:58:05"
$ws.Cells.Item(493, 4).Value = 33.8
$ws.Cells.Item(493, 5).Value = 21.3
$ws.Cells.Item(493, 6).Value = " "
$ws.Cells.Item(493, 7).Value = "   "

$ws.Cells.Item(494, 1).Value = 493
$ws.Cells.Item(494, 2).Value = "28.11.2023"
$ws.Cells.Item(494, 3).Value = "14:59:06"
$ws.Cells.Item(494, 4).Value = 33.8
$ws.Cells.Item(494, 5).Value = 21.2
$ws.Cells.Item(494, 6).Value = " "
$ws.Cells.Item(494, 7).Value = "   "

$ws.Cells.Item(495, 1).Value = 494
$ws.Cells.Item(495, 2).Value = "28.11.2023"
$ws.Cells.Item(495, 3).Value = "15:00:07"
$ws.Cells.Item(495, 4).Value = 33.8
$ws.Cells.Item(495, 5).Value = 21.2
$ws.Cells.Item(495, 6).Value = " "
$ws.Cells.Item(495, 7).Value = "   "

$ws.Cells.Item(496, 1).Value = 495
$ws.Cells.Item(496, 2).Value = "28.11.2023"
$ws.Cells.Item(496, 3).Value = "15:01:08"
$ws.Cells.Item(496, 4).Value = 34
$ws.Cells.Item(496, 5).Value = 21.1
$ws.Cells.Item(496, 6).Value = " "
$ws.Cells.Item(496, 7).Value = "   "

$ws.Cells.Item(497, 1).Value = 496
$ws.Cells.Item(497, 2).Value = "28.11.2023"
$ws.Cells.Item(497, 3).Value = "15:05:56"
$ws.Cells.Item(497, 4).Value = 34.8
$ws.Cells.Item(497, 5).Value = 21
$ws.Cells.Item(497, 6).Value = " "
$ws.Cells.Item(497, 7).Value = "   "

